# Apply a refreshed cryptocurrency price/volume scrape to the "cryptos"
# worksheet (columns D = Price, E = Volume(1h), rows 2-51), matching the
# commit "Updated cryptos list on Mon Jul 10 04:39:28 UTC 2023 with GitHub
# Actions".
#
# Price/percentage values are stored as literal text in the source sheet
# (e.g. "30.174.51", "9.330", "+0.55%") rather than numbers, so each cell is
# forced to a text number format before the value is written -- otherwise
# Excel would auto-coerce numeric-looking strings (dropping significant
# trailing zeros, or turning multi-dot strings into dates/errors). The
# number format / style is reset back to Normal immediately afterwards so
# the cell's formatting matches the original (unstyled) cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = '30.174.51'
    "D3" = '1.863.75'
    "E3" = '  -0.42%  '
    "E4" = '  -0.12%  '
    "D5" = '234.02'
    "E5" = '  -1.09%  '
    "E6" = '  -0.09%  '
    "D7" = '0.4671'
    "E7" = '  -0.90%  '
    "D8" = '43.21'
    "E8" = '  +0.55%  '
    "D9" = '0.2862'
    "E9" = '  -1.09%  '
    "D10" = '0.06484'
    "E10" = '  -2.30%  '
    "D11" = '21.20'
    "E11" = '  -2.13%  '
    "D12" = '0.07742'
    "E12" = '  -3.88%  '
    "D13" = '1.882.73'
    "E13" = '  +0.60%  '
    "D14" = '93.73'
    "E14" = '  -3.85%  '
    "D15" = '0.6836'
    "E15" = '  -0.74%  '
    "D16" = '5.043'
    "E16" = '  -1.96%  '
    "D17" = '268.69'
    "E17" = '  -1.31%  '
    "D18" = '30.162.65'
    "E18" = '  -0.54%  '
    "E19" = '  -6.19%  '
    "D20" = '0.000007624'
    "E20" = '  -1.64%  '
    "E21" = '  -0.05%  '
    "D22" = '2.122.01'
    "E22" = '  +0.21%  '
    "E23" = '  -0.06%  '
    "D24" = '5.148'
    "E24" = '  -3.34%  '
    "D25" = '6.112'
    "E25" = '  -1.82%  '
    "D26" = '9.330'
    "E26" = '  -0.24%  '
    "D27" = '165.39'
    "E27" = '  -1.66%  '
    "D28" = '18.56'
    "E28" = '  -2.34%  '
    "D29" = '1.893'
    "E29" = '  -3.33%  '
    "D30" = '1.362'
    "D31" = '0.09872'
    "E31" = '  -1.14%  '
    "D32" = '1.450'
    "E32" = '  -1.13%  '
    "D33" = '4.235'
    "E33" = '  -3.23%  '
    "D34" = '4.008'
    "E34" = '  -2.05%  '
    "D35" = '0.04670'
    "E35" = '  -0.91%  '
    "D36" = '1.117'
    "E36" = '  -1.60%  '
    "D37" = '0.6877'
    "E37" = '  -2.19%  '
    "D38" = '2.704'
    "E38" = '  -0.39%  '
    "D39" = '0.01833'
    "E39" = '  -2.84%  '
    "D40" = '2.751'
    "E40" = '  +3.71%  '
    "D41" = '6.307'
    "E41" = '  -0.33%  '
    "D42" = '71.35'
    "E42" = '  -2.05%  '
    "E43" = '  -0.02%  '
    "D44" = '1.893'
    "E44" = '  -3.69%  '
    "D45" = '0.8341'
    "E45" = '  -1.18%  '
    "D46" = '101.99'
    "E46" = '  -1.29%  '
    "D47" = '0.4055'
    "E47" = '  -2.85%  '
    "D48" = '935.13'
    "E48" = '  -0.33%  '
    "D49" = '9.088'
    "E49" = '  -1.83%  '
    "D50" = '6.962'
    "E50" = '  -2.15%  '
    "D51" = '33.96'
    "E51" = '  -1.86%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
